# Add a "File" type attribute ("xfile") to the TypeTest entity.
#
# This mirrors the author's workflow:
#  1. On the "org_molgenis_test_TypeTest" sheet, a new column (AS) is
#     appended to the header row for the new "xfile" attribute.
#  2. On the "attributes" metadata sheet, a new row describing the
#     "xfile" attribute (entity=org_molgenis_test_TypeTest,
#     dataType=file, refEntity=FileMeta, idAttribute=false,
#     nillable=true) is inserted right before the existing
#     "xcomputedxref" row, shifting every row below it down by one.
#  3. The workbook is left with the TypeTest sheet active/selected
#     (it was "attributes" before the edit).

$wb = $excel.ActiveWorkbook

$wsTypeTest   = $wb.Worksheets.Item("org_molgenis_test_TypeTest")
$wsAttributes = $wb.Worksheets.Item("attributes")

# --- 1. "attributes" sheet: insert the new "xfile" attribute row -----------
# Insert a fresh row at 49 (pushes xcomputedxref/xcomputedint/Chromosome/
# Position/id/age/driverslicence rows down to 50-56) and populate it.
$wsAttributes.Rows.Item(49).Insert()

$wsAttributes.Range("A49").Value = "xfile"
$wsAttributes.Range("B49").Value = "org_molgenis_test_TypeTest"
$wsAttributes.Range("C49").Value = "file"
$wsAttributes.Range("D49").Value = "FileMeta"
$wsAttributes.Range("E49").Value = $false
$wsAttributes.Range("F49").Value = $true

# --- 2. "org_molgenis_test_TypeTest" sheet: add the "xfile" header ---------
$wsTypeTest.Range("AS1").Value = "xfile"

# --- 3. Page setup on "attributes" (portrait, paper size 9/A4) ------------
$wsAttributes.PageSetup.PaperSize = 9
$wsAttributes.PageSetup.Orientation = 1

# --- 4. Leave the UI state pointed at the TypeTest sheet/new column -------
$wsAttributes.Range("T49").Select()

$wsTypeTest.Activate()
$wsTypeTest.Range("AW21").Select()
